# LinkedIn Visuals - findEmployeeCount: refresh sample employee-count data
# on the "Data" sheet and drop the now-unused trailing blank row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Helper: write a value into a cell as TEXT (shared string), not as a
# number, even when the text looks numeric (e.g. "15515"). A direct
# `.Value = "15515"` assignment would be auto-converted to a numeric
# cell by Excel's normal type inference, which also drags in a new
# number-format style. Routing the literal through a text formula and
# pasting-special as values keeps the original cell style untouched and
# keeps the cell typed as a shared string, matching how these sample
# figures were refreshed in the sheet.
function Set-TextValue([object]$range, [string]$text) {
    $helper = $ws.Range("Z100")
    $escaped = $text.Replace('"', '""')
    $helper.Formula = '="' + $escaped + '"'
    $helper.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# Row 3 - Raymond James
Set-TextValue $ws.Range("D3") "10001+"
Set-TextValue $ws.Range("E3") "15515"

# Row 4 - TresVista
Set-TextValue $ws.Range("D4") "1,001-5,000"
Set-TextValue $ws.Range("E4") "138"

# Row 5 - Goldman Sachs
Set-TextValue $ws.Range("D5") "10001+"
Set-TextValue $ws.Range("E5") "57474"

# Row 6 - Morgan Stanley
Set-TextValue $ws.Range("D6") "10001+"
Set-TextValue $ws.Range("E6") "86317"

$ws.Range("Z100").Clear()

# Column F (the blank trailing column) now matches the plain bordered
# look already used by columns A-C instead of its own one-off style.
$ws.Range("A3").Copy()
$ws.Range("F3:F6").PasteSpecial(-4122)  # xlPasteFormats

# The trailing blank row (row 7) is no longer part of the table.
$ws.Rows.Item(7).Delete()

# Move the selection to reflect where the editor left off.
$ws.Range("E12").Select()
